# The commit adds one new weekly price-report row for "Vega Modelo de
# Temuco" / Pomelo (Start Ruby, Primera) dated 2022-08-09, inserted right
# before the existing row for 2022-07-06 (currently row 226). Inserting a
# row there shifts all subsequent rows (226-275) down to (227-276), which
# matches the dimension growing from A1:T275 to A1:T276.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(226).Insert()

$ws.Cells.Item(226, 1).Value = 10
$ws.Cells.Item(226, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(226, 3).Value = 'La Araucanía'
$ws.Cells.Item(226, 4).Value = 44782
$ws.Cells.Item(226, 5).Value = 9
$ws.Cells.Item(226, 6).Value = 'Fruta'
$ws.Cells.Item(226, 7).Value = 100102
$ws.Cells.Item(226, 8).Value = 'Cítricos'
$ws.Cells.Item(226, 9).Value = 100102006
$ws.Cells.Item(226, 10).Value = 'Pomelo'
$ws.Cells.Item(226, 11).Value = 'Start Ruby'
$ws.Cells.Item(226, 12).Value = 'Primera'
$ws.Cells.Item(226, 13).Value = 55
$ws.Cells.Item(226, 14).Value = 10000
$ws.Cells.Item(226, 15).Value = 12000
$ws.Cells.Item(226, 16).Value = 10909
$ws.Cells.Item(226, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(226, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(226, 19).Value = 727
$ws.Cells.Item(226, 20).Value = 15
